$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.46"
$ws.Range("E2").Value = "'1.06%"
$ws.Range("D3").Value = "'31.64"
$ws.Range("E3").Value = "'1.92%"
$ws.Range("D4").Value = "'5.004"
$ws.Range("E4").Value = "'0.97%"
$ws.Range("D5").Value = "'0.07688"
$ws.Range("E5").Value = "'4.77%"
$ws.Range("D6").Value = "'2.258"
$ws.Range("E6").Value = "'-1.81%"
$ws.Range("D7").Value = "'7.888"
$ws.Range("E7").Value = "'2.10%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9250"
$ws.Range("E8").Value = "'1.96%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.09692"
$ws.Range("E9").Value = "'21.14%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1745"
$ws.Range("E10").Value = "'3.93%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08396"
$ws.Range("E11").Value = "'2.85%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03252"
$ws.Range("E12").Value = "'4.88%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09854"
$ws.Range("E13").Value = "'-2.32%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001472"
$ws.Range("E14").Value = "'-3.16%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005771"
$ws.Range("E15").Value = "'-0.75%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.502"
$ws.Range("E16").Value = "'0.38%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.789"
$ws.Range("E17").Value = "'1.27%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.142"
$ws.Range("E18").Value = "'3.29%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3362"
$ws.Range("E19").Value = "'0.94%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1324"
$ws.Range("E20").Value = "'1.57%"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "'4.045"
$ws.Range("E21").Value = "'1.68%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2278"
$ws.Range("E22").Value = "'8.67%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04500"
$ws.Range("E23").Value = "'-1.07%"
$ws.Range("D24").Value = "'0.001211"
$ws.Range("E24").Value = "'-0.10%"
$ws.Range("D25").Value = "'0.004354"
$ws.Range("E25").Value = "'-6.48%"
$ws.Range("D26").Value = "'0.0001285"
$ws.Range("E26").Value = "'-1.08%"
$ws.Range("D27").Value = "'0.0003368"
$ws.Range("E27").Value = "'-0.77%"
$ws.Range("D39").Value = "'0.01684"
$ws.Range("E39").Value = "'4.79%"
$ws.Range("D40").Value = "'0.04628"
$ws.Range("E40").Value = "'4.14%"
$ws.Range("D41").Value = "'0.007503"
$ws.Range("E41").Value = "'2.30%"
$ws.Range("D42").Value = "'0.009755"
$ws.Range("E42").Value = "'10.66%"
$ws.Range("D43").Value = "'0.1384"
$ws.Range("E43").Value = "'4.15%"
$ws.Range("D44").Value = "'0.002125"
$ws.Range("E44").Value = "'6.27%"
$ws.Range("D45").Value = "'0.009403"
$ws.Range("E45").Value = "'-1.15%"
$ws.Range("D46").Value = "'0.00006068"
$ws.Range("E46").Value = "'2.35%"
$ws.Range("D47").Value = "'0.00000000744"
$ws.Range("E47").Value = "'-0.70%"
$ws.Range("D48").Value = "'2.794"
$ws.Range("E48").Value = "'24.69%"
$ws.Range("D49").Value = "'0.001984"
$ws.Range("E49").Value = "'-31.52%"
$ws.Range("D50").Value = "'0.00002084"
$ws.Range("E50").Value = "'-0.70%"
$ws.Range("D51").Value = "'0.0001985"
$ws.Range("E51").Value = "'-0.70%"
